$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.602.16'
$ws.Range('E2').Value = '  -2.37%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.841.95'
$ws.Range('E3').Value = '  -1.37%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  -0.11%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.35'
$ws.Range('E5').Value = '  -1.51%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9999'
$ws.Range('E6').Value = '  -0.17%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4247'
$ws.Range('E7').Value = '  -2.72%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3641'
$ws.Range('E8').Value = '  -1.62%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.82'
$ws.Range('E9').Value = '  +1.86%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07264'
$ws.Range('E10').Value = '  -3.36%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.8978'
$ws.Range('E11').Value = '  -4.54%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '20.61'
$ws.Range('E12').Value = '  -3.96%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.813.73'
$ws.Range('E13').Value = '  -4.85%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.380'
$ws.Range('E14').Value = '  -1.44%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.555'
$ws.Range('E15').Value = '  -2.53%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.06862'
$ws.Range('E16').Value = '  -0.06%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.002'
$ws.Range('E17').Value = '  -0.09%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '78.22'
$ws.Range('E18').Value = '  -5.04%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000008853'
$ws.Range('E19').Value = '  -2.67%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9981'
$ws.Range('E20').Value = '  -0.41%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.57'
$ws.Range('E21').Value = '  -2.61%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '27.606.55'
$ws.Range('E22').Value = '  -2.35%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.968'
$ws.Range('E23').Value = '  -3.25%  '

$ws.Range('E24').Value = '  -1.75%  '

$ws.Range('B25').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C25').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.016.58'
$ws.Range('E25').Value = '  -5.38%  '

$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.043'
$ws.Range('E26').Value = '  +0.92%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '154.10'
$ws.Range('E27').Value = '  -0.53%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.25'
$ws.Range('E28').Value = '  -1.01%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.245'
$ws.Range('E29').Value = '  -1.57%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.832'
$ws.Range('E30').Value = '  +5.85%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '110.70'
$ws.Range('E31').Value = '  -2.75%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.08864'
$ws.Range('E32').Value = '  -1.86%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7756'
$ws.Range('E33').Value = '  -3.24%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.557'
$ws.Range('E34').Value = '  -5.80%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.938'
$ws.Range('E35').Value = '  -0.44%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.092'
$ws.Range('E36').Value = '  -6.80%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9983'
$ws.Range('E37').Value = '  -0.28%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05431'
$ws.Range('E38').Value = '  -0.13%  '

$ws.Range('E39').Value = '  -2.11%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01923'
$ws.Range('E40').Value = '  -1.83%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.809'
$ws.Range('E41').Value = '  -5.55%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5067'
$ws.Range('E42').Value = '  -3.51%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.795'
$ws.Range('E43').Value = '  -4.83%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1639'
$ws.Range('E44').Value = '  -2.18%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.214'
$ws.Range('E45').Value = '  -5.85%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.06637'
$ws.Range('E46').Value = '  -1.97%  '

$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.36'
$ws.Range('E47').Value = '  -1.33%  '

$ws.Range('B48').Value = 'Decentraland'
$ws.Range('C48').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4718'
$ws.Range('E48').Value = '  -3.31%  '

$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '105.62'
$ws.Range('E49').Value = '  -2.31%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.9992'
$ws.Range('E50').Value = '  -0.11%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.636'
$ws.Range('E51').Value = '  -2.63%  '
